$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text representation (avoid Excel
# auto-converting numeric-looking strings like "27.263.05" or "0.07300"
# into floating point numbers, which would lose formatting/precision).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.263.05"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.857.11"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.85%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.07"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.29%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4659"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.71%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3709"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07300"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8911"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07870"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.855.82"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.412"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.517"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008933"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.12%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.299.93"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.25%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.091.89"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.029"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +9.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.68"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.41"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.047"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.99"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.048"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08847"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.145"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7693"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +5.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.169"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.528"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.728"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +10.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.107"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.25%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.41%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.945"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.075"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1627"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.529"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4798"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.59%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.99"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.648"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.81%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.54"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.65%  "
